$d = $word.ActiveDocument

$replacements = @(
    @{Old = "96×61=5856"; New = "97×46=4462"},
    @{Old = "29×21=609";  New = "39×77=3003"},
    @{Old = "57×82=4674"; New = "76×59=4484"},
    @{Old = "38×18=684";  New = "45×41=1845"},
    @{Old = "48×48=2304"; New = "45×38=1710"},
    @{Old = "90×41=3690"; New = "11×46=506"},
    @{Old = "48×93=4464"; New = "82×91=7462"},
    @{Old = "97×73=7081"; New = "93×92=8556"},
    @{Old = "83×63=5229"; New = "73×31=2263"},
    @{Old = "72×93=6696"; New = "20×17=340"},
    @{Old = "36×18=648";  New = "32×95=3040"},
    @{Old = "26×98=2548"; New = "24×77=1848"},
    @{Old = "37×96=3552"; New = "60×78=4680"},
    @{Old = "19×46=874";  New = "70×56=3920"},
    @{Old = "42×90=3780"; New = "36×33=1188"},
    @{Old = "25×71=1775"; New = "33×92=3036"},
    @{Old = "80×23=1840"; New = "79×52=4108"},
    @{Old = "32×88=2816"; New = "62×24=1488"},
    @{Old = "67×25=1675"; New = "73×62=4526"},
    @{Old = "36×56=2016"; New = "47×71=3337"},
    @{Old = "80×79=6320"; New = "18×30=540"},
    @{Old = "50×27=1350"; New = "38×79=3002"},
    @{Old = "74×85=6290"; New = "96×23=2208"},
    @{Old = "79×51=4029"; New = "87×32=2784"},
    @{Old = "21×45=945";  New = "93×53=4929"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
